$d = $word.ActiveDocument

# Collapse the stray run of extra spaces between "Romano," and "Bill Harrison"
# (originally spread across three separate runs, one of them wrapped in
# <w:proofErr> gramStart/gramEnd tags) down to a single space, matching
# the corrected text "Maddalena Romano, Bill Harrison".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Maddalena Romano,          Bill Harrison", $false, $false, $false, $false, $false, $true, 1, $false, "Maddalena Romano, Bill Harrison", 2)
